$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "('climat', 'chang')"
$ws.Range("C2").Value = 171
$ws.Range("B3").Value = "('global', 'warm')"
$ws.Range("C3").Value = 36
$ws.Range("B4").Value = "('sea', 'level')"
$ws.Range("C4").Value = 20
$ws.Range("B5").Value = "('polit', 'statement')"
$ws.Range("C5").Value = 16
$ws.Range("B6").Value = "('al', 'gore')"
$ws.Range("C6").Value = 16
$ws.Range("B7").Value = "('year', 'ago')"
$ws.Range("C7").Value = 14
$ws.Range("B8").Value = "('3rd', 'world')"
$ws.Range("C8").Value = 14
$ws.Range("B9").Value = "('video', 'game')"
$ws.Range("C9").Value = 13
$ws.Range("B10").Value = "('level', 'rise')"
$ws.Range("C10").Value = 13
$ws.Range("B11").Value = "('reddit', 'kotakuinact')"
$ws.Range("C11").Value = 13
$ws.Range("B12").Value = "('kotakuinact', 'comment')"
$ws.Range("C12").Value = 13
$ws.Range("B13").Value = "('ice', 'cap')"
$ws.Range("C13").Value = 12
$ws.Range("B14").Value = "('last', 'year')"
$ws.Range("C14").Value = 12
$ws.Range("B15").Value = "('man', 'make')"
$ws.Range("C15").Value = 11
$ws.Range("B16").Value = "('ice', 'age')"
$ws.Range("C16").Value = 9
$ws.Range("B17").Value = "('hockey', 'stick')"
$ws.Range("C17").Value = 9
$ws.Range("B18").Value = "('late', 'game')"
$ws.Range("C18").Value = 9
$ws.Range("B19").Value = "('year', 'year')"
$ws.Range("C19").Value = 9
$ws.Range("B20").Value = "('gather', 'storm')"
$ws.Range("C20").Value = 8
$ws.Range("B21").Value = "('co2', 'emiss')"
$ws.Range("C21").Value = 8
$ws.Range("B22").Value = "('carbon', 'emiss')"
$ws.Range("C22").Value = 8
$ws.Range("B23").Value = "('nasa', 'gov')"
$ws.Range("C23").Value = 8
$ws.Range("B24").Value = "('pari', 'agreement')"
$ws.Range("C24").Value = 8
$ws.Range("B25").Value = "('specif', 'heat')"
$ws.Range("C25").Value = 8
$ws.Range("B26").Value = "('power', 'plant')"
$ws.Range("C26").Value = 7
$ws.Range("B27").Value = "('peopl', 'think')"
$ws.Range("C27").Value = 7
$ws.Range("B28").Value = "('co2', 'level')"
$ws.Range("C28").Value = 7
$ws.Range("B29").Value = "('black', 'peopl')"
$ws.Range("C29").Value = 7
$ws.Range("B30").Value = "('low', 'iq')"
$ws.Range("C30").Value = 7
$ws.Range("B31").Value = "('coal', 'oil')"
$ws.Range("C31").Value = 7
$ws.Range("B32").Value = "('polit', 'issu')"
$ws.Range("C32").Value = 6
$ws.Range("B33").Value = "('melt', 'ice')"
$ws.Range("C33").Value = 6
$ws.Range("B34").Value = "('effect', 'climat')"
$ws.Range("C34").Value = 6
$ws.Range("B35").Value = "('carbon', 'dioxid')"
$ws.Range("C35").Value = 6
$ws.Range("B36").Value = "('realli', 'want')"
$ws.Range("C36").Value = 6
$ws.Range("B37").Value = "('lord', 'believ')"
$ws.Range("C37").Value = 6
$ws.Range("B38").Value = "('chang', 'polit')"
$ws.Range("C38").Value = 6
$ws.Range("B39").Value = "('degre', 'celsius')"
$ws.Range("C39").Value = 6
$ws.Range("B40").Value = "('climat', 'scientist')"
$ws.Range("C40").Value = 6
$ws.Range("B41").Value = "('climat', 'scienc')"
$ws.Range("C41").Value = 6
$ws.Range("B42").Value = "('chang', 'real')"
$ws.Range("C42").Value = 6
$ws.Range("B43").Value = "('global', 'climat')"
$ws.Range("C43").Value = 6
$ws.Range("B44").Value = "('greenhous', 'effect')"
$ws.Range("C44").Value = 6
$ws.Range("B45").Value = "('bell', 'curv')"
$ws.Range("C45").Value = 6
$ws.Range("B46").Value = "('tile', 'flood')"
$ws.Range("C46").Value = 6
$ws.Range("B47").Value = "('publish', 'report')"
$ws.Range("C47").Value = 5
$ws.Range("B48").Value = "('carbon', 'pollut')"
$ws.Range("C48").Value = 5
$ws.Range("B49").Value = "('global', 'catastroph')"
$ws.Range("C49").Value = 5
$ws.Range("B50").Value = "('thing', 'happen')"
$ws.Range("C50").Value = 5
$ws.Range("B51").Value = "('mind', 'worm')"
$ws.Range("C51").Value = 5
$ws.Range("B52").Value = "('climat', 'model')"
$ws.Range("C52").Value = 5
$ws.Range("B53").Value = "('natur', 'disast')"
$ws.Range("C53").Value = 5
$ws.Range("B54").Value = "('fossil', 'fuel')"
$ws.Range("C54").Value = 5
$ws.Range("B55").Value = "('chang', 'thing')"
$ws.Range("C55").Value = 5
$ws.Range("B56").Value = "('say', 'thing')"
$ws.Range("C56").Value = 5
$ws.Range("B57").Value = "('settl', 'scienc')"
$ws.Range("C57").Value = 5
$ws.Range("B58").Value = "('make', 'sen')"
$ws.Range("C58").Value = 5
$ws.Range("B59").Value = "('human', 'hive')"
$ws.Range("C59").Value = 5
$ws.Range("B60").Value = "('giss', 'nasa')"
$ws.Range("C60").Value = 5
$ws.Range("B61").Value = "('believ', 'climat')"
$ws.Range("C61").Value = 5
$ws.Range("B62").Value = "('go', 'back')"
$ws.Range("C62").Value = 5
$ws.Range("B63").Value = "('chang', 'happen')"
$ws.Range("C63").Value = 5
$ws.Range("B64").Value = "('high', 'co2')"
$ws.Range("C64").Value = 5
$ws.Range("B65").Value = "('feel', 'good')"
$ws.Range("C65").Value = 5
$ws.Range("B66").Value = "('green', 'tech')"
$ws.Range("C66").Value = 5
$ws.Range("B67").Value = "('solar', 'panel')"
$ws.Range("C67").Value = 5
$ws.Range("B68").Value = "('long', 'time')"
$ws.Range("C68").Value = 5
$ws.Range("B69").Value = "('chang', 'climat')"
$ws.Range("C69").Value = 5
$ws.Range("B70").Value = "('polar', 'bear')"
$ws.Range("C70").Value = 5
$ws.Range("B71").Value = "('grow', 'wheat')"
$ws.Range("C71").Value = 5
$ws.Range("B72").Value = "('get', 'grip')"
$ws.Range("C72").Value = 5
$ws.Range("B73").Value = "('say', 'year')"
$ws.Range("C73").Value = 5
$ws.Range("B74").Value = "('take', 'much')"
$ws.Range("C74").Value = 5
$ws.Range("B75").Value = "('use', 'coal')"
$ws.Range("C75").Value = 5
$ws.Range("B76").Value = "('make', 'think')"
$ws.Range("C76").Value = 5
$ws.Range("B77").Value = "('think', 'peopl')"
$ws.Range("C77").Value = 5
$ws.Range("B78").Value = "('import', 'peopl')"
$ws.Range("C78").Value = 5
$ws.Range("B79").Value = "('hundr', 'year')"
$ws.Range("C79").Value = 5
$ws.Range("B80").Value = "('corn', 'seed')"
$ws.Range("C80").Value = 5
$ws.Range("B81").Value = "('unit', 'nation')"
$ws.Range("C81").Value = 4
